$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order matters: it determines the order new shared strings are appended in,
# which must match the authored workbook's edit history.
$ws.Range("F5").Value = "Whole Team"
$ws.Range("E5").Value = " NA"
$ws.Range("D4").Value = "Yes"
$ws.Range("F4").Value = "Caitlin / Sav"
$ws.Range("F7").Value = "Caitlin/Sav"
$ws.Range("E4").Value = "NA"

$ws.Range("D5").Value = "Yes"
$ws.Range("D6").Value = "Yes"
$ws.Range("D7").Value = "Yes"

$ws.Range("E6").Value = "NA"
$ws.Range("E7").Value = "NA"

$ws.Range("F6").Value = "Caitlin / Sav"

$ws.Range("G7").Value = 2

# Caitlin's task moved from "Group Evaluation" to "Sprint Chart Generator";
# the now-obsolete "Group Evaluation" row (22) is cleared out.
$ws.Range("B21").Value = "Sprint Chart Generator"

$ws.Range("A22").Value = ""
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("G22").Value = ""
